$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
$rng = $ws.Range("ZZ500")
$rng.Font.Name = "Tahoma"
$rng.Font.Size = 11
$rng.Font.Color = 4089918
$ws.Rows.Item(500).Delete()
